$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at row 439 (everything currently at/after row 439
# shifts down by two), mirroring the weekly update for
# "Feria Lagunitas de Puerto Montt - Repollo".
$ws.Rows.Item(439).Insert()
$ws.Rows.Item(439).Insert()

# --- Row 439: Copenhague / Primera ---
$ws.Cells.Item(439, 1).Value = 4
$ws.Cells.Item(439, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(439, 3).Value = "Los Lagos"
$ws.Cells.Item(439, 4).Value = 44746
$ws.Cells.Item(439, 5).Value = 10
$ws.Cells.Item(439, 6).Value = 100112006
$ws.Cells.Item(439, 7).Value = "Repollo"
$ws.Cells.Item(439, 8).Value = "Copenhague"
$ws.Cells.Item(439, 9).Value = "Primera"
$ws.Cells.Item(439, 10).Value = 500
$ws.Cells.Item(439, 11).Value = 2000
$ws.Cells.Item(439, 12).Value = 2000
$ws.Cells.Item(439, 13).Value = 2000
$ws.Cells.Item(439, 14).Value = "`$/unidad"
$ws.Cells.Item(439, 15).Value = "Región Metropolitana"
$ws.Cells.Item(439, 16).Value = 2000
$ws.Cells.Item(439, 17).Value = 1
$ws.Cells.Item(439, 18).Value = "Hortaliza"

# --- Row 440: Crespo record / Primera ---
$ws.Cells.Item(440, 1).Value = 4
$ws.Cells.Item(440, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(440, 3).Value = "Los Lagos"
$ws.Cells.Item(440, 4).Value = 44746
$ws.Cells.Item(440, 5).Value = 10
$ws.Cells.Item(440, 6).Value = 100112006
$ws.Cells.Item(440, 7).Value = "Repollo"
$ws.Cells.Item(440, 8).Value = "Crespo record"
$ws.Cells.Item(440, 9).Value = "Primera"
$ws.Cells.Item(440, 10).Value = 500
$ws.Cells.Item(440, 11).Value = 1800
$ws.Cells.Item(440, 12).Value = 1900
$ws.Cells.Item(440, 13).Value = 1850
$ws.Cells.Item(440, 14).Value = "`$/unidad"
$ws.Cells.Item(440, 15).Value = "Región del Maule"
$ws.Cells.Item(440, 16).Value = 1850
$ws.Cells.Item(440, 17).Value = 1
$ws.Cells.Item(440, 18).Value = "Hortaliza"
